# reporteDiarioCobranza.docx edit
#
# 1. Grow the data row (row 2) from trHeight 215 -> 1140 twips.
# 2. The "numeroPago" cell:
#      "{numeroPago}-{#credito}{numerosPagos}" / "{/credito}"
#    becomes
#      "{numeroPago}- {numeroPagos}" / "" (the credito loop markers are
#      removed, the second paragraph is kept but left empty).
# 3. The "interes" cell gains a leading "$" so it reads "${interes}"
#    like the neighbouring pago/iva/seguro/cargo cells.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(2)

# 1. Row height: 215 twips -> 1140 twips (Word COM Height is in points;
#    1 point = 20 twips).
$row.Height = 1140 / 20

# 2a. Collapse "{numeroPago}-{#credito}{numerosPagos}" to
#     "{numeroPago}- {numeroPagos}".
$d.Content.Find.Execute('{numeroPago}-{#credito}{numerosPagos}', $true, $false, $false, $false, $false, $true, 1, $false, '{numeroPago}- {numeroPagos}', 2) | Out-Null

# 2b. Empty out the paragraph that used to hold the "{/credito}" closing tag
#     (the paragraph itself stays, just with no text left in it).
$d.Content.Find.Execute('{/credito}', $true, $false, $false, $false, $false, $true, 1, $false, '', 2) | Out-Null

# 3. "{interes}" -> "${interes}"
$d.Content.Find.Execute('{interes}', $true, $false, $false, $false, $false, $true, 1, $false, '${interes}', 2) | Out-Null
